$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.143.87"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.573.55"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Formula = "'517.91"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("D6").Formula = "'142.05"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D7").Formula = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Formula = "'0.565"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").Value = "2.588.76"
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").Formula = "'6.75"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("E12").Value = "  -4.07%  "
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "3.026.42"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").Value = "58.106.67"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Formula = "'20.32"
$ws.Range("E16").Value = "  -2.73%  "
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "2.571.52"
$ws.Range("E18").Value = "  -2.94%  "
$ws.Range("D19").Formula = "'341.77"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Formula = "'65.62"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("E26").Value = "  -5.61%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "2.686.84"
$ws.Range("E28").Value = "  -2.77%  "
$ws.Range("D29").Formula = "'6.98"
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("D30").Value = "0.0₃0744"
$ws.Range("E30").Value = "  -6.72%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Formula = "'6.23"
$ws.Range("E32").Value = "  -6.70%  "
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("D34").Formula = "'18.70"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("D35").Formula = "'149.85"
$ws.Range("E35").Value = "  -1.80%  "
$ws.Range("D36").Formula = "'4.00"
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("D38").Formula = "'0.866"
$ws.Range("E38").Value = "  -5.29%  "
$ws.Range("D39").Formula = "'35.96"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("D40").Formula = "'0.833"
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("D41").Formula = "'1.44"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("E42").Value = "  -3.31%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Formula = "'269.97"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Formula = "'10.66"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Formula = "'0.0950"
$ws.Range("E46").Value = "  -1.92%  "
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("D48").Formula = "'18.84"
$ws.Range("E48").Value = "  -2.93%  "
$ws.Range("D49").Formula = "'0.0522"
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("D50").Value = "1.970.93"
$ws.Range("E50").Value = "  -3.44%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Formula = "'18.39"
$ws.Range("E51").Value = "  +0.27%  "
